# Add a third "UPPRESSO_PPTOKEN_MODE" time-test data block in columns Q:W,
# mirroring the existing MITREDID (A:G) / UPPRESSO_RING_MODE (I:O) blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- column widths for the new block (Q..V) -------------------------------
# (values tuned so the engine's char-width quantization lands on the pixel
# width closest to the target stored width)
$ws.Columns.Item(17).ColumnWidth = 25.15   # Q -> ~25.89
$ws.Columns.Item(18).ColumnWidth = 15.72   # R -> ~16.44
$ws.Columns.Item(19).ColumnWidth = 10.85   # S -> ~11.55
$ws.Columns.Item(20).ColumnWidth = 13.85   # T -> ~14.55
$ws.Columns.Item(21).ColumnWidth = 12.15   # U -> ~12.89
$ws.Columns.Item(22).ColumnWidth = 10.00   # V -> ~10.66

# ---- header row (row 2) ---------------------------------------------------
$headers = @("UPPRESSO_PPTOKEN_MODE","Request_Generate","Token_Verify","Request_Process","Code_Process","Total_Time","rtt")
$col = 17
foreach ($h in $headers) {
    $ws.Cells.Item(2, $col).Value = $h
    $col++
}

# ---- data rows (rows 3-12): run#, Request_Generate, Token_Verify, Request_Process, Code_Process, Total_Time
$data = @(
    @(1, 808, 18, 1020, 378, 6005),
    @(2, 808, 19, 1114, 372, 6255),
    @(3, 786, 23, 1171, 407, 5934),
    @(4, 802, 23, 1171, 358, 6063),
    @(5, 879, 19, 1231, 363, 6412),
    @(6, 791, 25,  965, 372, 5581),
    @(7, 874, 19, 1196, 382, 6168),
    @(8, 736, 21, 1244, 375, 6245),
    @(9, 813, 20, 1079, 405, 6011),
    @(10, 937, 22, 1030, 378, 6249)
)

$row = 3
foreach ($vals in $data) {
    $ws.Cells.Item($row, 17).Value = $vals[0]   # Q - run index
    $ws.Cells.Item($row, 18).Value = $vals[1]   # R - Request_Generate
    $ws.Cells.Item($row, 19).Value = $vals[2]   # S - Token_Verify
    $ws.Cells.Item($row, 20).Value = $vals[3]   # T - Request_Process
    $ws.Cells.Item($row, 21).Value = $vals[4]   # U - Code_Process
    $ws.Cells.Item($row, 22).Value = $vals[5]   # V - Total_Time
    $ws.Cells.Item($row, 23).Value = 0          # W - rtt
    $row++
}

# ---- AVG row (row 13) ------------------------------------------------------
$ws.Range("Q13").Value = "AVG"
$ws.Range("R13").Formula = "=AVERAGE(R3:R12)"
$ws.Range("S13").Formula = "=AVERAGE(S3:S12)"
$ws.Range("T13").Formula = "=AVERAGE(T3:T12)"
$ws.Range("U13").Formula = "=AVERAGE(U3:U12)"
$ws.Range("V13").Formula = "=AVERAGE(V3:V12)"
$ws.Range("W13").Value = 0

# ---- formatting: center horizontally + vertically, like the I:O block -----
$ws.Range("Q2:W13").HorizontalAlignment = -4108   # xlCenter
$ws.Range("Q2:W13").VerticalAlignment = -4108     # xlCenter

# ---- move the view / selection to roughly match the edited area -----------
$excel.ActiveWindow.ScrollColumn = 11
[void]$ws.Range("T20").Select()
